# Types doorgevoerd in hpm en ich en excel bestanden
#
# HPM.waarde("...") calls are split into HPM.waardeNum("...") / HPM.waardeString("...")
# depending on whether the underlying HPM field is numeric or textual, on the three
# "Interpretatie" header sheets. Frozen-pane/selection view state and autofit column
# widths are refreshed to reflect the new (longer) header labels.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: InterpretatieCF -------------------------------------------------
$ws1 = $wb.Worksheets.Item("InterpretatieCF")

$ws1.Range("A1").Value = 'HPM.waardeNum("IRT")'
$ws1.Range("B1").Value = 'HPM.waardeNum("PAP")'
$ws1.Range("C1").Value = 'HPM.waardeNum("DNA_N")'
$ws1.Range("D1").Value = 'HPM.waardeString("DNA_imut1")'
$ws1.Range("E1").Value = 'HPM.waardeString("DNA_imut2")'

$ws1.Activate()
$ws1.Range("A2").Select()

# --- Sheet 2: InterpretatieEGABlad1 ------------------------------------------
$ws2 = $wb.Worksheets.Item("InterpretatieEGABlad1")

$ws2.Range("A1").Value = 'HPM.waardeNum("EGA_N")'
$ws2.Range("B1").Value = 'HPM.waardeString("EGA_imut1")'
$ws2.Range("C1").Value = 'HPM.waardeString("EGA_imut2")'

$ws2.Columns.Item(1).ColumnWidth = 28.33
$ws2.Columns.Item(2).ColumnWidth = 33
$ws2.Columns.Item(3).ColumnWidth = 33

$ws2.Activate()
$ws2.Range("A1:D1048576").Select()

# --- Sheet 3: InterpretatieEGABlad2 ------------------------------------------
$ws3 = $wb.Worksheets.Item("InterpretatieEGABlad2")

$ws3.Range("A1").Value = 'HPM.waardeString("EGA_imut1")'
$ws3.Range("B1").Value = 'HPM.waardeString("EGA_imut2")'

$ws3.Columns.Item(1).ColumnWidth = 33
$ws3.Columns.Item(2).ColumnWidth = 33

$ws3.Activate()
$ws3.Range("A1:C1048576").Select()

# --- Restore the originally active/selected tab -------------------------------
$ws1.Activate()
